# Cost Modeling First Attempt
# Developed basic cost modeling test case in the ISP reviews workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.28515625
$ws.Columns.Item(5).ColumnWidth = 11.5703125
$ws.Range("J1:L1").EntireColumn.ColumnWidth = 12.42578125

# ---------------------------------------------------------------------
# Row 50 - section title (bold)
# ---------------------------------------------------------------------
$ws.Range("A50").Value = "NASA Advanced Mission Cost Model (from HSMAD)"
$ws.Range("A50").Font.Bold = $true

# ---------------------------------------------------------------------
# Row 51 - alpha
# ---------------------------------------------------------------------
$ws.Range("A51").Value = "alpha"
$ws.Range("B51").Value = [double]"5.6499999999999996E-4"
$ws.Range("B51").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------
# Row 52 - beta
# ---------------------------------------------------------------------
$ws.Range("A52").Value = "beta"
$ws.Range("B52").Value = 0.5941

# ---------------------------------------------------------------------
# Row 53 - zi
# ---------------------------------------------------------------------
$ws.Range("A53").Value = "zi"
$ws.Range("B53").Value = 0.6604

# ---------------------------------------------------------------------
# Row 54 - delta, inflation hyperlink, inflation multiplier
# ---------------------------------------------------------------------
$ws.Range("A54").Value = "delta"
$ws.Range("B54").Value = 80.599
$ws.Range("E54").Value = "1999-2016 Inflation"
$ws.Hyperlinks.Add($ws.Range("E54"), "https://www.bls.gov/data/inflation_calculator.htm") | Out-Null
$ws.Range("G54").Value = 1.42

# ---------------------------------------------------------------------
# Row 55 - epsilon
# ---------------------------------------------------------------------
$ws.Range("A55").Value = "epsilon"
$ws.Range("B55").Value = [double]"3.8084999999999997E-55"
$ws.Range("B55").NumberFormat = "0.0000E+00"

# ---------------------------------------------------------------------
# Row 56 - phi, category headers (merged, bold, centered)
# ---------------------------------------------------------------------
$ws.Range("A56").Value = "phi"
$ws.Range("B56").Value = -0.3553

$ws.Range("E56").Value = "LH2 Engines"
$ws.Range("E56:H56").Merge()
$ws.Range("E56:H56").Font.Bold = $true
$ws.Range("E56:H56").HorizontalAlignment = -4108

$ws.Range("I56").Font.Bold = $true

$ws.Range("J56").Value = "Nuclear Thermal Rockets"
$ws.Range("J56:L56").Merge()
$ws.Range("J56:L56").Font.Bold = $true
$ws.Range("J56:L56").HorizontalAlignment = -4108

$ws.Range("M56").Font.Bold = $true

$ws.Range("N56").Value = "Electric Propulsion"
$ws.Range("N56:Q56").Merge()
$ws.Range("N56:Q56").Font.Bold = $true
$ws.Range("N56:Q56").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Row 57 - gamma, Isp headers for each engine column (bold italic)
# ---------------------------------------------------------------------
$ws.Range("A57").Value = "gamma"
$ws.Range("B57").Value = 1.5691

$ws.Range("E57").Value = 445
$ws.Range("F57").Value = 452
$ws.Range("G57").Value = 465
$ws.Range("H57").Value = 480
$ws.Range("E57:H57").Font.Bold = $true
$ws.Range("E57:H57").Font.Italic = $true
$ws.Range("E57:H57").HorizontalAlignment = -4108

$ws.Range("J57").Value = 850
$ws.Range("K57").Value = 950
$ws.Range("L57").Value = 1000
$ws.Range("N57").Value = 3000
$ws.Range("O57").Value = 3800
$ws.Range("P57").Value = 5000
$ws.Range("Q57").Value = 9000
$ws.Range("I57:Q57").Font.Bold = $true
$ws.Range("I57:Q57").Font.Italic = $true

# ---------------------------------------------------------------------
# Row 58 - Q / quantity
# ---------------------------------------------------------------------
$ws.Range("A58").Value = "Q"
$ws.Range("B58").Value = "quantity"
$ws.Range("E58").Value = 1
$ws.Range("F58").Value = 1
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1
$ws.Range("J58").Value = 1
$ws.Range("K58").Value = 1
$ws.Range("L58").Value = 1

# ---------------------------------------------------------------------
# Row 59 - M / dry mass (lbs)  (formulas, 0-decimal number format)
# ---------------------------------------------------------------------
$ws.Range("A59").Value = "M"
$ws.Range("B59").Value = "dry mass (lbs)"
$ws.Range("E59").Formula = '=F59*0.95'
$ws.Range("F59").Formula = '=F69*2.2/F58'
$ws.Range("G59").Formula = '=F59*1.1'
$ws.Range("H59").Formula = '=F59*1.25'
$ws.Range("J59").Formula = '=J69*2.2'
$ws.Range("K59").Formula = '=J59*1.1'
$ws.Range("L59").Formula = '=J59*1.15'
$ws.Range("E59:L59").NumberFormat = "0"

# ---------------------------------------------------------------------
# Row 60 - S / specification
# ---------------------------------------------------------------------
$ws.Range("A60").Value = "S"
$ws.Range("B60").Value = "specification"
$ws.Range("E60").Value = 2.39
$ws.Range("F60").Value = 2.39
$ws.Range("G60").Value = 2.39
$ws.Range("H60").Value = 2.39
$ws.Range("J60").Value = 2.39
$ws.Range("K60").Value = 2.39
$ws.Range("L60").Value = 2.39
$ws.Range("N60").Value = 2.39
$ws.Range("O60").Value = 2.39
$ws.Range("P60").Value = 2.39
$ws.Range("Q60").Value = 2.39

# ---------------------------------------------------------------------
# Row 61 - IOC / initial operating capability
# ---------------------------------------------------------------------
$ws.Range("A61").Value = "IOC"
$ws.Range("B61").Value = "initial operating capability"
$ws.Range("E61").Value = 2030
$ws.Range("F61").Value = 2030
$ws.Range("G61").Value = 2030
$ws.Range("H61").Value = 2030
$ws.Range("J61").Value = 2030
$ws.Range("K61").Value = 2030
$ws.Range("L61").Value = 2030
$ws.Range("N61").Value = 2030
$ws.Range("O61").Value = 2030
$ws.Range("P61").Value = 2030
$ws.Range("Q61").Value = 2030

# ---------------------------------------------------------------------
# Row 62 - B / Block number
# ---------------------------------------------------------------------
$ws.Range("A62").Value = "B "
$ws.Range("B62").Value = "Block number"
$ws.Range("E62").Value = 3
$ws.Range("F62").Value = 3
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 1
$ws.Range("J62").Value = 1
$ws.Range("K62").Value = 1
$ws.Range("L62").Value = 1
$ws.Range("N62").Value = 1
$ws.Range("O62").Value = 1
$ws.Range("P62").Value = 1
$ws.Range("Q62").Value = 1

# ---------------------------------------------------------------------
# Row 63 - D / Difficulty
# ---------------------------------------------------------------------
$ws.Range("A63").Value = "D"
$ws.Range("B63").Value = "Difficulty"
$ws.Range("E63").Value = -1.5
$ws.Range("F63").Value = -0.5
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 1
$ws.Range("J63").Value = 1.5
$ws.Range("K63").Value = 1.5
$ws.Range("L63").Value = 2
$ws.Range("N63").Value = 0.5
$ws.Range("O63").Value = 0.5
$ws.Range("P63").Value = 1
$ws.Range("Q63").Value = 1.5

# ---------------------------------------------------------------------
# Row 66 - Static Engine Mass (kg) + reviewer comment
# ---------------------------------------------------------------------
$ws.Range("A66").Value = "Static Engine Mass (kg)"
$ws.Range("E66").Value = 0
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 34500
$ws.Range("J66").NumberFormat = "#,##0"
$ws.Range("K66").Formula = '=J66*1.05'
$ws.Range("L66").Formula = '=J66*1.1'
$ws.Range("N66").Value = 0
$ws.Range("O66").Value = 0
$ws.Range("P66").Value = 0
$ws.Range("Q66").Value = 0
$ws.Range("R66").Value = "What about the electric propulsion nuclear generator or solar panel mass?  NTR reactor mass seems high?"

# ---------------------------------------------------------------------
# Row 67 - Inert Mass Ratio
# ---------------------------------------------------------------------
$ws.Range("A67").Value = "Inertr Mass Ratio"
$ws.Range("E67").Value = 0.17
$ws.Range("F67").Value = 0.18
$ws.Range("G67").Value = 0.19
$ws.Range("H67").Value = 0.2
$ws.Range("J67").Value = 0.1
$ws.Range("K67").Value = 0.11
$ws.Range("L67").Value = 0.12
$ws.Range("N67").Value = 0.18
$ws.Range("O67").Value = 0.19
$ws.Range("P67").Value = 0.2
$ws.Range("Q67").Value = 0.23

# ---------------------------------------------------------------------
# Row 68 - Prop Mass (kg)
# ---------------------------------------------------------------------
$ws.Range("A68").Value = "Prop Mass (kg)"
$ws.Range("F68").Value = 125060
$ws.Range("J68").Value = 92847

# ---------------------------------------------------------------------
# Row 69 - Engine Mass (kg)
# ---------------------------------------------------------------------
$ws.Range("A69").Value = "Engine Mass (kg)"
$ws.Range("F69").Value = 21260
$ws.Range("J69").Formula = '=J66+J68*0.1'
$ws.Range("J69").NumberFormat = "0"

# ---------------------------------------------------------------------
# Row 70 - totals
# ---------------------------------------------------------------------
$ws.Range("F70").Formula = '=SUM(F68:F69)'
$ws.Range("J70").Formula = '=SUM(J68:J69)'

# ---------------------------------------------------------------------
# Row 73 - Total Development Cost ($M 1999) + reviewer comment
# ---------------------------------------------------------------------
$ws.Range("A73").Value = "Total Development Cost (`$M 1999)"
$ws.Range("E73").Formula = '=$B$51*E58^$B$52*E59*$B$53*$B$54^E60*$B$55^(1/(E61-1900))*E62^$B$56*$B$57^E63'
$ws.Range("F73").Formula = '=$B$51*F58^$B$52*F59*$B$53*$B$54^F60*$B$55^(1/(F61-1900))*F62^$B$56*$B$57^F63'
$ws.Range("G73").Formula = '=$B$51*G58^$B$52*G59*$B$53*$B$54^G60*$B$55^(1/(G61-1900))*G62^$B$56*$B$57^G63'
$ws.Range("H73").Formula = '=$B$51*H58^$B$52*H59*$B$53*$B$54^H60*$B$55^(1/(H61-1900))*H62^$B$56*$B$57^H63'
$ws.Range("J73").Formula = '=$B$51*J58^$B$52*J59*$B$53*$B$54^J60*$B$55^(1/(J61-1900))*J62^$B$56*$B$57^J63'
$ws.Range("K73").Formula = '=$B$51*K58^$B$52*K59*$B$53*$B$54^K60*$B$55^(1/(K61-1900))*K62^$B$56*$B$57^K63'
$ws.Range("L73").Formula = '=$B$51*L58^$B$52*L59*$B$53*$B$54^L60*$B$55^(1/(L61-1900))*L62^$B$56*$B$57^L63'
$ws.Range("R73").Value = "These prices are WAY too high`u{2026} the mass seems to be really large for these based on what I'd expect"
$ws.Range("E73:R73").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

# ---------------------------------------------------------------------
# Row 74 - Total Development Cost ($M 2016)
# ---------------------------------------------------------------------
$ws.Range("A74").Value = "Total Development Cost (`$M 2016)"
$ws.Range("E74").Formula = '=E73*$G$54'
$ws.Range("F74").Formula = '=F73*$G$54'
$ws.Range("G74").Formula = '=G73*$G$54'
$ws.Range("H74").Formula = '=H73*$G$54'
$ws.Range("J74").Formula = '=J73*$G$54'
$ws.Range("K74").Formula = '=K73*$G$54'
$ws.Range("L74").Formula = '=L73*$G$54'
$ws.Range("E74:R74").NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

# ---------------------------------------------------------------------
# Row 77 - note
# ---------------------------------------------------------------------
$ws.Range("A77").Value = "SSME weights 7,775"

# ---------------------------------------------------------------------
# Rows 79-82 - source references
# ---------------------------------------------------------------------
$ws.Range("A79").Value = "http://www.nasa.gov/sites/default/files/files/01_CEH_Main_Body_02_27_15.pdf"
$ws.Range("A80").Value = "http://www.nasa.gov/sites/default/files/files/CEH_AppD.pdf"
$ws.Range("A81").Value = "http://ntrs.nasa.gov/archive/nasa/casi.ntrs.nasa.gov/20140005476.pdf"
$ws.Range("A82").Value = "http://ntrs.nasa.gov/archive/nasa/casi.ntrs.nasa.gov/20140005340.pdf"

# ---------------------------------------------------------------------
# Sheet view / page setup
# ---------------------------------------------------------------------
$ws.Range("F58").Select()
$excel.ActiveWindow.ScrollRow = 43
$ws.PageSetup.Orientation = 1
